$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original "General" style while forcing text storage:
# set NumberFormat to Text ("@") right before assigning each value so
# Excel keeps the literal string instead of coercing it to a number/percentage.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.42"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.60%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.715"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.18%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06215"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.22%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.728"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.82%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.39%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9123"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.65%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1401"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.15%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.81%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07088"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.24%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03111"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.61%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09055"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.13%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001529"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.41%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006160"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.22%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005975"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.66%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.450"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.02%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.179"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.18%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.167"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.87%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1310"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.68%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.095"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.08%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04227"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.30%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001180"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.49%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004074"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.15%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.01%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.11%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03946"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.01%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1112"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.02%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004144"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.88%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.38%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01315"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-19.62%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.55%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.03%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-37.52%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2485"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "82.96%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
